# Update countries & provincias Spain
# Applies the data refresh for the "Pais" sheet: updates the "last updated"
# timestamp, refreshes case counters for several countries, and re-sorts
# Jordania/Burkina Faso (Jordania's total now exceeds Burkina Faso's, so it
# moves one row up in the (descending, by total cases) table).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 19:22"

# --- Estados Unidos (row 4) --------------------------------------------
$ws.Range("B4").Value = 352160
$ws.Range("C4").Value = 15487
$ws.Range("E4").Value = 322544
$ws.Range("F4").Value = 8832
$ws.Range("G4").Value = 753
$ws.Range("H4").Value = 10369

# --- Francia (row 8) ----------------------------------------------------
$ws.Range("B8").Value = 98010
$ws.Range("C8").Value = 5171
$ws.Range("D8").Value = 17250
$ws.Range("E8").Value = 71849
$ws.Range("F8").Value = 7072
$ws.Range("G8").Value = 833
$ws.Range("H8").Value = 8911

# --- Turquia (row 12) ----------------------------------------------------
$ws.Range("F12").Value = 1415

# --- Suiza (row 13) -------------------------------------------------------
$ws.Range("D13").Value = 8056
$ws.Range("E13").Value = 12834

# --- Austria (row 17) ------------------------------------------------------
$ws.Range("B17").Value = 12280
$ws.Range("C17").Value = 229
$ws.Range("E17").Value = 8597

# --- Irlanda (row 26) -------------------------------------------------------
$ws.Range("B26").Value = 5364
$ws.Range("C26").Value = 370
$ws.Range("E26").Value = 5165
$ws.Range("G26").Value = 16
$ws.Range("H26").Value = 174

# --- Kazajistan (row 76) ------------------------------------------------------
$ws.Range("D76").Value = 46
$ws.Range("E76").Value = 599

# --- Afganistan (row 91) --------------------------------------------------------
$ws.Range("D91").Value = 18
$ws.Range("E91").Value = 342

# --- Jordania / Burkina Faso reorder (rows 94-95) -------------------------
# Jordania's total cases rose to 349, overtaking Burkina Faso (345), so the
# two rows swap places; Burkina Faso's figures are carried down unchanged.
$ws.Range("A94").Value = "Jordania"
$ws.Range("B94").Value = 349
$ws.Range("C94").Value = 4
$ws.Range("D94").Value = 126
$ws.Range("E94").Value = 217
$ws.Range("F94").Value = 5
$ws.Range("G94").Value = 1
$ws.Range("H94").Value = 6

$ws.Range("A95").Value = "Burkina Faso"
$ws.Range("B95").Value = 345
$ws.Range("C95").Value = 0
$ws.Range("D95").Value = 90
$ws.Range("E95").Value = 238
$ws.Range("F95").Value = 0
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 17

# --- Liechtenstein (row 129) -----------------------------------------------------
$ws.Range("D129").Value = 55
$ws.Range("E129").Value = 21

# --- Zambia (row 148) --------------------------------------------------------------
$ws.Range("D148").Value = 5
$ws.Range("E148").Value = 33
